# GitHub Actions daily update: append the latest gold-price row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The data starts at row 1 (header) and grows by one row per run;
# find the first empty row right after the current last used row.
$newRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($newRow, 1).Value = "31-10-2025"
$ws.Cells.Item($newRow, 2).Value = "The price of gold in India today is ₹12,268 per gram for 24 karat gold, ₹11,245 per gram for 22 karat gold and ₹9,201 per gram for 18 karat gold (also called 999 gold)."

# Match the formatting already used by the existing data rows
# (thin border on both cells, word-wrap on the long text column).
$ws.Cells.Item($newRow, 1).Borders.LineStyle = 1
$ws.Cells.Item($newRow, 2).Borders.LineStyle = 1
$ws.Cells.Item($newRow, 2).WrapText = $true
